$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Content update: swap the two "Key Personnel" rows.
#   Row 5 (person #1) becomes Mason Tandal / BICSI ITS Technician
#   Row 6 (person #2) becomes David Yokoi / BICSI RCDD (moved up from row 5)
#   Ryan Yokoi's entry is removed entirely.
# ---------------------------------------------------------------------------
$ws.Range("B5").Value = "Tandal, Mason M."
$ws.Range("E5").Value = "BICSI ITS Technician"
$ws.Range("F5").Value = "808-722-7257"
$ws.Range("G5").Value = "mason.tandal@bcshawaii.com"

$ws.Range("B6").Value = "Yokoi, David R."
$ws.Range("E6").Value = "BICSI RCDD"
$ws.Range("F6").Value = "808-585-7999"
$ws.Range("G6").Value = "david@bcshawaii.com"

# ---------------------------------------------------------------------------
# Formatting: thicken the interior vertical dividers of the header row
# (between Certification / Contact # / E-mail Address) from thin to medium.
# ---------------------------------------------------------------------------
$xlEdgeLeft = 7
$xlEdgeRight = 10
$xlContinuous = 1
$xlLineStyleNone = -4142
$xlMedium = -4138

$ws.Range("E4").Borders.Item($xlEdgeLeft).LineStyle = $xlContinuous
$ws.Range("E4").Borders.Item($xlEdgeLeft).Weight = $xlMedium
$ws.Range("E4").Borders.Item($xlEdgeRight).LineStyle = $xlLineStyleNone

$ws.Range("F4").Borders.Item($xlEdgeLeft).LineStyle = $xlContinuous
$ws.Range("F4").Borders.Item($xlEdgeLeft).Weight = $xlMedium
$ws.Range("F4").Borders.Item($xlEdgeRight).LineStyle = $xlContinuous
$ws.Range("F4").Borders.Item($xlEdgeRight).Weight = $xlMedium

$ws.Range("G4").Borders.Item($xlEdgeLeft).LineStyle = $xlLineStyleNone

# ---------------------------------------------------------------------------
# Column G is widened (manual resize instead of the previous best-fit width).
# ---------------------------------------------------------------------------
$ws.Columns.Item(7).ColumnWidth = 27.93

# ---------------------------------------------------------------------------
# Move the active selection to G10.
# ---------------------------------------------------------------------------
$ws.Range("G10").Select() | Out-Null
